$d = $word.ActiveDocument

# Disable "smart quotes" autocorrect so literal straight apostrophes we
# type below are not turned into curly quotes.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

# ---------------------------------------------------------------------
# Change 1 (paragraph 1): the folder path at the end of the sentence is
# replaced with a new path pointing at the GitHub repo location.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    ": /Users/rudolfo/LightFieldMicroscopy/Simulation/Birefringence/2024_02/",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found1) {
    $rng1.Text = ": '/Users/rudolfo/Software/GitHub/BirTomo/data/2025_02/SpiculeA Experim&Simulation/Simulation Data'"
}

# ---------------------------------------------------------------------
# Change 2: note that "further processed" now explicitly mentions
# rectification as one of the processing steps.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "further processed.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)
if ($found2) {
    $rng2.Text = "further processed, including rectification."
}
